# Commiting scripts (R22 UAT2 - Regression).
#
# Appends 40 new LCY Saving Accounts records (rows 72-111) to Sheet0,
# mirroring the additional TC / Account_ID / Customer_ID / PD rows that
# were committed to "Data/LCY Saving Accounts.xlsx".
#
# All four columns hold digit-only identifiers that must stay TEXT (the
# existing rows 1-71 store every value as a shared string, never as a
# number). Assigning a numeric-looking string straight to .Value/.Value2
# makes Excel auto-convert it to a number, so each value is written with
# a leading apostrophe (forces text entry) and the quote-prefix/"Text"
# style that trick stamps onto the cell is stripped back off with
# ClearFormats() afterwards - leaving plain text cells identical in kind
# to the pre-existing data, with no residual formatting differences.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rows = @(
  @("118518","1008784314","17898803","6004"),
  @("118518","1008784315","17898804","6005"),
  @("118518","1008784336","17898816","6004"),
  @("118518","1008784337","17898817","6005"),
  @("118518","1008784351","17898848","6004"),
  @("118518","1008784360","17898858","6004"),
  @("118518","1008784361","17898859","6005"),
  @("118518","1008784370","17898868","6004"),
  @("118518","1008784378","17898879","6004"),
  @("118518","1008784379","17898881","6005"),
  @("118518","1008784387","17898893","6004"),
  @("118518","1008784388","17898894","6005"),
  @("118518","1008784395","17898908","6004"),
  @("118518","1008784396","17898909","6005"),
  @("118500","1008784400","17898936","6004"),
  @("118518","1008784405","17898940","6005"),
  @("118500","1008784408","17898948","6004"),
  @("118500","1008784411","17898951","6004"),
  @("118500","1008784415","17898954","6004"),
  @("118500","1008784419","17898967","6004"),
  @("118518","1008784422","17898970","6005"),
  @("118500","1008784425","17898992","6004"),
  @("118518","1008784428","17898996","6005"),
  @("118500","1008784430","17899005","6004"),
  @("118518","1008784433","17899008","6005"),
  @("118500","1008784434","17899013","6004"),
  @("118518","1008784437","17899016","6005"),
  @("118500","1008784444","17899025","6004"),
  @("118518","1008784447","17899028","6005"),
  @("118518","1008784450","17899032","6005"),
  @("118518","1008784453","17899035","6004"),
  @("118518","1008784455","17899037","6005"),
  @("118518","1008784459","17899043","6005"),
  @("118518","1008784462","17899046","6004"),
  @("118518","1008784468","17899055","6005"),
  @("118518","1008784469","17899058","6004"),
  @("118518","1008784472","17899061","6005"),
  @("118518","1008784477","17899067","6005"),
  @("118518","1008784479","17899070","6004"),
  @("118518","1008784482","17899073","6005"),
)

$startRow = 72
$endRow = $startRow + $rows.Count - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = "'" + $vals[$c]
    }
}

# Drop the quote-prefix ("Text") style the apostrophe entries above applied -
# the source workbook keeps these shared-string cells on the default style.
$ws.Range("A" + $startRow + ":D" + $endRow).ClearFormats()

